$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly refresh of Fruta/Hortaliza dataset: rotate Fecha, Volumen, Precio
# minimo/maximo/promedio ponderado, Origen and Precio $/Kg across rows 2-16
# per the updated source extract (diff maps old row values onto new rows).

# Row 2
$ws.Range("D2").Value2 = 44533
$ws.Range("M2").Value2 = 150
$ws.Range("N2").Value2 = 4000
$ws.Range("O2").Value2 = 4000
$ws.Range("P2").Value2 = 4000
$ws.Range("R2").Value = "Provincia de Curicó"
$ws.Range("S2").Value2 = 2000

# Row 3
$ws.Range("D3").Value2 = 44978
$ws.Range("M3").Value2 = 500
$ws.Range("N3").Value2 = 3000
$ws.Range("O3").Value2 = 3000
$ws.Range("P3").Value2 = 3000
$ws.Range("R3").Value = "Provincia de Curicó"
$ws.Range("S3").Value2 = 1500

# Row 4
$ws.Range("D4").Value2 = 44194
$ws.Range("M4").Value2 = 120
$ws.Range("N4").Value2 = 3000
$ws.Range("O4").Value2 = 3000
$ws.Range("P4").Value2 = 3000
$ws.Range("R4").Value = "Provincia de Linares"
$ws.Range("S4").Value2 = 1500

# Row 5
$ws.Range("D5").Value2 = 44174
$ws.Range("M5").Value2 = 200
$ws.Range("N5").Value2 = 3200
$ws.Range("O5").Value2 = 3200
$ws.Range("P5").Value2 = 3200
$ws.Range("R5").Value = "Provincia de Curicó"
$ws.Range("S5").Value2 = 1600

# Row 6
$ws.Range("D6").Value2 = 44208
$ws.Range("M6").Value2 = 85
$ws.Range("N6").Value2 = 3000
$ws.Range("O6").Value2 = 3000
$ws.Range("P6").Value2 = 3000
$ws.Range("R6").Value = "Provincia de Linares"
$ws.Range("S6").Value2 = 1500

# Row 7
$ws.Range("D7").Value2 = 44586
$ws.Range("M7").Value2 = 250
$ws.Range("N7").Value2 = 5000
$ws.Range("O7").Value2 = 5000
$ws.Range("P7").Value2 = 5000
$ws.Range("R7").Value = "Provincia de Curicó"
$ws.Range("S7").Value2 = 2500

# Row 8
$ws.Range("D8").Value2 = 44582
$ws.Range("M8").Value2 = 380
$ws.Range("N8").Value2 = 5000
$ws.Range("O8").Value2 = 5000
$ws.Range("P8").Value2 = 5000
$ws.Range("R8").Value = "Provincia de Curicó"
$ws.Range("S8").Value2 = 2500

# Row 9
$ws.Range("D9").Value2 = 44236
$ws.Range("M9").Value2 = 300
$ws.Range("N9").Value2 = 3600
$ws.Range("O9").Value2 = 4000
$ws.Range("P9").Value2 = 3800
$ws.Range("R9").Value = "Provincia de Curicó"
$ws.Range("S9").Value2 = 1900

# Row 10
$ws.Range("D10").Value2 = 44980
$ws.Range("M10").Value2 = 250
$ws.Range("N10").Value2 = 4000
$ws.Range("O10").Value2 = 4000
$ws.Range("P10").Value2 = 4000
$ws.Range("R10").Value = "Provincia de Curicó"
$ws.Range("S10").Value2 = 2000

# Row 11
$ws.Range("D11").Value2 = 44232
$ws.Range("M11").Value2 = 200
$ws.Range("N11").Value2 = 3000
$ws.Range("O11").Value2 = 3000
$ws.Range("P11").Value2 = 3000
$ws.Range("R11").Value = "Provincia de Curicó"
$ws.Range("S11").Value2 = 1500

# Row 12
$ws.Range("D12").Value2 = 44188
$ws.Range("M12").Value2 = 150
$ws.Range("N12").Value2 = 3000
$ws.Range("O12").Value2 = 3400
$ws.Range("P12").Value2 = 3240
$ws.Range("R12").Value = "Provincia de Linares"
$ws.Range("S12").Value2 = 1620

# Row 13
$ws.Range("D13").Value2 = 44238
$ws.Range("M13").Value2 = 300
$ws.Range("N13").Value2 = 3600
$ws.Range("O13").Value2 = 4000
$ws.Range("P13").Value2 = 3800
$ws.Range("R13").Value = "Provincia de Curicó"
$ws.Range("S13").Value2 = 1900

# Row 14
$ws.Range("D14").Value2 = 44168
$ws.Range("M14").Value2 = 170
$ws.Range("N14").Value2 = 8000
$ws.Range("O14").Value2 = 8000
$ws.Range("P14").Value2 = 8000
$ws.Range("R14").Value = "Provincia de Linares"
$ws.Range("S14").Value2 = 4000

# Row 15
$ws.Range("D15").Value2 = 44617
$ws.Range("M15").Value2 = 90
$ws.Range("N15").Value2 = 6500
$ws.Range("O15").Value2 = 6500
$ws.Range("P15").Value2 = 6500
$ws.Range("R15").Value = "Provincia de Curicó"
$ws.Range("S15").Value2 = 3250

# Row 16
$ws.Range("D16").Value2 = 44231
$ws.Range("M16").Value2 = 150
$ws.Range("N16").Value2 = 3400
$ws.Range("O16").Value2 = 3400
$ws.Range("P16").Value2 = 3400
$ws.Range("R16").Value = "Provincia de Curicó"
$ws.Range("S16").Value2 = 1700
